$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated category (x) values in column B and updated count (y) values in
# column C for rows 3-22, which feed the bar chart's cached series data.
$newData = @(
    @{ Row = 3;  B = -3.3880685093911098;   C = 36 },
    @{ Row = 4;  B = -3.0806869660071499;   C = 838 },
    @{ Row = 5;  B = -2.7733054226231899;   C = 5024 },
    @{ Row = 6;  B = -2.46592387923923;     C = 11580 },
    @{ Row = 7;  B = -2.1585423358552802;   C = 23766 },
    @{ Row = 8;  B = -1.85116079247132;     C = 44450 },
    @{ Row = 9;  B = -1.5437792490873601;   C = 52212 },
    @{ Row = 10; B = -1.2363977057033999;   C = 56276 },
    @{ Row = 11; B = -0.92901616231944595;  C = 61278 },
    @{ Row = 12; B = -0.62163461893548899;  C = 67183 },
    @{ Row = 13; B = -0.31425307555153098;  C = 62896 },
    @{ Row = 14; B = -0.0068715321675740801;C = 52902 },
    @{ Row = 15; B = 0.30051001121638299;   C = 52667 },
    @{ Row = 16; B = 0.60789155460034106;   C = 51614 },
    @{ Row = 17; B = 0.91527309798429801;   C = 35512 },
    @{ Row = 18; B = 1.22265464136826;      C = 25722 },
    @{ Row = 19; B = 1.5300361847522099;    C = 16699 },
    @{ Row = 20; B = 1.8374177281361701;    C = 4677 },
    @{ Row = 21; B = 2.1447992715201298;    C = 561 },
    @{ Row = 22; B = 2.4521808149040898;    C = 25 }
)

foreach ($item in $newData) {
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 3).Value = $item.C
}

# Restore the selection that was active when the workbook was last saved.
$ws.Range("A1:N23").Select()
